$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H88").Value = 1598.2
$ws.Range("I88").Value = 1303
$ws.Range("J88").Value = 1672
$ws.Range("K88").Value = 1303
$ws.Range("L88").Value = 1672
$ws.Range("M88").Value = -897
$ws.Range("N88").Value = -2484
$ws.Range("H91").Value = 1598.2
$ws.Range("I91").Value = 1303
$ws.Range("J91").Value = 1672
$ws.Range("K91").Value = 1303
$ws.Range("L91").Value = 1672
$ws.Range("M91").Value = 101
$ws.Range("N91").Value = -4480
$ws.Range("H103").Value = 9091527
$ws.Range("I103").Value = 439.45456
$ws.Range("K103").Value = 1318.36368
$ws.Range("M103").Value = -732.3636799999999
$ws.Range("H116").Value = 5147.1665
$ws.Range("I116").Value = 3775.9
$ws.Range("K116").Value = 3775.9
$ws.Range("M116").Value = -333.9000000000001
$ws.Range("H137").Value = 1907.4324
$ws.Range("I137").Value = 943
$ws.Range("J137").Value = 3916.6667
$ws.Range("K137").Value = 2829
$ws.Range("L137").Value = 11750.0001
$ws.Range("M137").Value = -279
$ws.Range("N137").Value = -16850.0001
$ws.Range("H138").Value = 2112.7856
$ws.Range("I138").Value = 1628.0416
$ws.Range("K138").Value = 4884.1248
$ws.Range("M138").Value = 255.8752000000004

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 6315.38
$ws.Range("I32").Value = 3444.115
$ws.Range("J32").Value = 25530.77
$ws.Range("K32").Value = 3444.115
$ws.Range("L32").Value = 25530.77
$ws.Range("M32").Value = -3157.115
$ws.Range("N32").Value = -26104.77
$ws.Range("H97").Value = 614.5
$ws.Range("I97").Value = 111.3
$ws.Range("J97").Value = 1872.5
$ws.Range("K97").Value = 111.3
$ws.Range("L97").Value = 1872.5
$ws.Range("M97").Value = 384.7
$ws.Range("N97").Value = -2864.5

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620
$ws.Range("H82").Value = 32192.875
$ws.Range("I82").Value = 14900
$ws.Range("K82").Value = 14900
$ws.Range("M82").Value = -14517
$ws.Range("H85").Value = 32192.875
$ws.Range("I85").Value = 14900
$ws.Range("K85").Value = 14900
$ws.Range("M85").Value = -13574
$ws.Range("H122").Value = 48500
$ws.Range("J122").Value = 48500
$ws.Range("L122").Value = 48500
$ws.Range("N122").Value = -58300
$ws.Range("H124").Value = 42265
$ws.Range("J124").Value = 42265
$ws.Range("L124").Value = 42265
$ws.Range("N124").Value = -52085
$ws.Range("H125").Value = 50374
$ws.Range("J125").Value = 50374
$ws.Range("L125").Value = 50374
$ws.Range("N125").Value = -60214

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H41").Value = 20409.285
$ws.Range("J41").Value = 20409.285
$ws.Range("L41").Value = 20409.285
$ws.Range("N41").Value = -21265.285
$ws.Range("H50").Value = 9093.799999999999
$ws.Range("J50").Value = 9093.799999999999
$ws.Range("L50").Value = 9093.799999999999
$ws.Range("N50").Value = -10343.8
$ws.Range("H51").Value = 9285.857
$ws.Range("J51").Value = 9285.857
$ws.Range("L51").Value = 9285.857
$ws.Range("N51").Value = -10757.857
$ws.Range("H60").Value = 30504.889
$ws.Range("J60").Value = 30504.889
$ws.Range("L60").Value = 30504.889
$ws.Range("N60").Value = -31526.889
$ws.Range("H61").Value = 9285.857
$ws.Range("J61").Value = 9285.857
$ws.Range("L61").Value = 9285.857
$ws.Range("N61").Value = -9981.857
$ws.Range("H68").Value = 16903.6
$ws.Range("J68").Value = 16903.6
$ws.Range("L68").Value = 16903.6
$ws.Range("N68").Value = -18401.6
$ws.Range("H71").Value = 16903.6
$ws.Range("J71").Value = 16903.6
$ws.Range("L71").Value = 50710.8
$ws.Range("N71").Value = -58198.8
$ws.Range("H99").Value = 2107.4285
$ws.Range("I99").Value = 1923.6
$ws.Range("J99").Value = 2274.5454
$ws.Range("K99").Value = 1923.6
$ws.Range("L99").Value = 2274.5454
$ws.Range("M99").Value = -425.5999999999999
$ws.Range("N99").Value = -5270.5454
$ws.Range("H109").Value = 10998
$ws.Range("J109").Value = 10998
$ws.Range("L109").Value = 10998
$ws.Range("N109").Value = -13078
$ws.Range("H126").Value = 2107.4285
$ws.Range("I126").Value = 1923.6
$ws.Range("J126").Value = 2274.5454
$ws.Range("K126").Value = 5770.799999999999
$ws.Range("L126").Value = 6823.6362
$ws.Range("M126").Value = -3300.799999999999
$ws.Range("N126").Value = -11763.6362
$ws.Range("H132").Value = 2086.7368
$ws.Range("I132").Value = 2281.2
$ws.Range("J132").Value = 1870.6666
$ws.Range("K132").Value = 6843.599999999999
$ws.Range("L132").Value = 5611.9998
$ws.Range("M132").Value = -4313.599999999999
$ws.Range("N132").Value = -10671.9998

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H22").Value = 925.125
$ws.Range("J22").Value = 1066.6666
$ws.Range("L22").Value = 3199.9998
$ws.Range("N22").Value = -3537.9998
$ws.Range("H27").Value = 925.125
$ws.Range("J27").Value = 1066.6666
$ws.Range("L27").Value = 3199.9998
$ws.Range("N27").Value = -3403.9998
$ws.Range("H44").Value = 603.50793
$ws.Range("I44").Value = 300
$ws.Range("J44").Value = 608.4032
$ws.Range("K44").Value = 900
$ws.Range("L44").Value = 1825.2096
$ws.Range("M44").Value = -502
$ws.Range("N44").Value = -2621.2096
$ws.Range("H119").Value = 2122.2693
$ws.Range("I119").Value = 1638.95
$ws.Range("J119").Value = 3733.3333
$ws.Range("K119").Value = 4916.85
$ws.Range("L119").Value = 11199.9999
$ws.Range("M119").Value = -78.85000000000036
$ws.Range("N119").Value = -20875.9999

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H41").Value = 3378.15
$ws.Range("I41").Value = 3500
$ws.Range("J41").Value = 3364.611
$ws.Range("K41").Value = 3500
$ws.Range("L41").Value = 3364.611
$ws.Range("M41").Value = -3145
$ws.Range("N41").Value = -4074.611
$ws.Range("H57").Value = 12558
$ws.Range("J57").Value = 20061
$ws.Range("L57").Value = 20061
$ws.Range("N57").Value = -21701
$ws.Range("H123").Value = 34315.75
$ws.Range("J123").Value = 34315.75
$ws.Range("L123").Value = 34315.75
$ws.Range("N123").Value = -39215.75
$ws.Range("H128").Value = 52620
$ws.Range("J128").Value = 52620
$ws.Range("L128").Value = 52620
$ws.Range("N128").Value = -62580
$ws.Range("H130").Value = 29996
$ws.Range("J130").Value = 29996
$ws.Range("L130").Value = 29996
$ws.Range("N130").Value = -40036
$ws.Range("H135").Value = 46338.094
$ws.Range("J135").Value = 46338.094
$ws.Range("L135").Value = 46338.094
$ws.Range("N135").Value = -56478.094

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H108").Value = 15996.667
$ws.Range("J108").Value = 15996.667
$ws.Range("L108").Value = 15996.667
$ws.Range("N108").Value = -23676.667

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H64").Value = 34556
$ws.Range("J64").Value = 34556
$ws.Range("L64").Value = 34556
$ws.Range("N64").Value = -35052
$ws.Range("H67").Value = 34556
$ws.Range("J67").Value = 34556
$ws.Range("L67").Value = 34556
$ws.Range("N67").Value = -36272
$ws.Range("H108").Value = 24690
$ws.Range("J108").Value = 24690
$ws.Range("L108").Value = 24690
$ws.Range("N108").Value = -32370
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("N109").ClearContents()
